$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Names"
$ws.Range("B1").Value = "Prices"

# Give the new header cell (B1) the same formatting as A1 (bold, bordered, centered)
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$data = @(
    @('The Legend of Zelda: Ocarina of Time', '91,99 €'),
    @('Super Mario Galaxy', '91,99 €'),
    @('Super Mario Galaxy 2', '91,99 €'),
    @('Metroid Prime', '89,99 €'),
    @('Super Mario Odyssey', '89,99 €'),
    @('Halo: Combat Evolved', '87,99 €'),
    @('The House in Fata Morgana - Dreams of the Revenants Edition -', '83,99 €'),
    @('NFL 2K1', '62,99 €'),
    @('Uncharted 2: Among Thieves', '88,99 €'),
    @('Tekken 3', '91,99 €'),
    @('The Legend of Zelda: The Wind Waker', '90,99 €'),
    @('Gran Turismo', '86,99 €'),
    @('Metal Gear Solid 2: Sons of Liberty', '88,99 €'),
    @('Grand Theft Auto Double Pack', '81,99 €'),
    @('Baldur''s Gate II: Shadows of Amn', '91,99 €'),
    @('Tetris Effect: Connected', '88,99 €'),
    @('The Legend of Zelda Collector''s Edition', '89,99 €'),
    @('Gran Turismo 3: A-Spec', '84,99 €'),
    @('The Legend of Zelda: A Link to the Past', '90,99 €'),
    @('The Legend of Zelda: Majora''s Mask', '91,99 €'),
    @('The Last of Us', '92,99 €'),
    @('Persona 5 Royal', '84,99 €'),
    @('The Last of Us Remastered', '92,99 €'),
    @('The Legend of Zelda: Ocarina of Time 3D', '90,99 €'),
    @('Chrono Cross', '88,99 €'),
    @('Gears of War', '84,99 €'),
    @('Sid Meier''s Civilization II', '88,99 €'),
    @('Halo 3', '81,99 €'),
    @('Ninja Gaiden Black', '88,99 €'),
    @('Super Mario Advance 4: Super Mario Bros. 3', '89,99 €'),
    @('Jet Grind Radio', '83,99 €'),
    @('Grim Fandango', '91,99 €'),
)
$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}
